$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections in existing rows ---
$ws.Range("A3").Value = "browserProcess"
$ws.Range("A9").Value = "fuelTypes"
$ws.Range("B12").Value = "Manual gearbox"

# --- New rows appended after the existing data ---
$ws.Range("A15").Value = "filePath"
$ws.Range("B15").Value = "Data\Output\OutputReport"

$ws.Range("A17").Value = "emailReceiver"
$ws.Range("B17").Value = "patricia.ciortin@fwfcompany.com;calin.gandila@fwfcompany.com;diana.sacacian@fwfcompany.com"

$ws.Range("A18").Value = "emailSubject"
$ws.Range("B18").Value = "Car List"

# --- Selection matches the saved workbook state ---
$ws.Range("A15").Select()
